$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.166.56"
$ws.Range("E2").Value = "  +0.96%  "

$ws.Range("D3").Value = "1.780.16"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.88%  "

$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.71"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.18%  "

$ws.Range("E9").Value = "  +1.15%  "

$ws.Range("E10").Value = "  +2.29%  "

$ws.Range("E11").Value = "  +1.03%  "

$ws.Range("D12").Value = "2.037.58"
$ws.Range("E12").Value = "  +0.29%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.58%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.775.94"
$ws.Range("E14").Value = "  -0.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.623"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.49%  "

$ws.Range("D16").Value = "34.140.66"
$ws.Range("E16").Value = "  +0.80%  "

$ws.Range("E17").Value = "  +1.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.78%  "

$ws.Range("E19").Value = "  +4.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "245.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.92%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.29%  "

$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("E23").Value = "  +2.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.46%  "

$ws.Range("E28").Value = "  +2.02%  "

$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("E30").Value = "  +2.54%  "

$ws.Range("E31").Value = "  +0.70%  "

$ws.Range("E32").Value = "  +4.27%  "

$ws.Range("E33").Value = "  +5.31%  "

$ws.Range("E34").Value = "  -0.99%  "

$ws.Range("D35").Value = "1.441.38"
$ws.Range("E35").Value = "  +3.90%  "

$ws.Range("E36").Value = "  +4.02%  "

$ws.Range("E37").Value = "  +7.11%  "

$ws.Range("E38").Value = "  +3.04%  "

$ws.Range("E39").Value = "  +0.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.66%  "

$ws.Range("E42").Value = "  +1.50%  "

$ws.Range("E43").Value = "  +0.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.63%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0511"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.28%  "

$ws.Range("E46").Value = "  +3.84%  "

$ws.Range("E47").Value = "  -0.05%  "

$ws.Range("D48").Value = "0.0₆0133"
$ws.Range("E48").Value = "  -2.90%  "

$ws.Range("D49").Value = "1.939.59"
$ws.Range("E49").Value = "  +0.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.87%  "

$ws.Range("E51").Value = "  +0.13%  "

